# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps for the
# zh-cn and de-de report sheets, rows 2 and 3.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 10:15:33"
$wsZhCn.Range("D3").Value = "2016-01-26 10:15:33"
$wsZhCn.Range("G2").Value = "2016-01-26 10:16:28"
$wsZhCn.Range("G3").Value = "2016-01-26 10:16:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 10:15:43"
$wsDeDe.Range("D3").Value = "2016-01-26 10:15:43"
$wsDeDe.Range("G2").Value = "2016-01-26 10:16:48"
$wsDeDe.Range("G3").Value = "2016-01-26 10:16:48"
